# Generate Report for Handoff
#
# The localization-status workbook tracks, per target language sheet
# (zh-cn / de-de), when each source file was last handed off for
# translation. Row 5 corresponds to file
#   a814842c-b9fd-44d1-8fa8-95ce81a1e82f.md
# which has just been (re-)handed off, so its "Latest Handoff Datetime"
# (column D) needs to be refreshed on both language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-03 08:43:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-03 08:43:38"
